$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the data up: drop the 2000/2005 rows, keep 2010/2015, and add a new 2020 row.
$ws.Range("A2").Value = "2010年"
$ws.Range("B2").Value = 77.37
$ws.Range("C2").Value = 74.83
$ws.Range("D2").Value = 72.38

$ws.Range("A3").Value = "2015年"
$ws.Range("B3").Value = 79.43000000000001
$ws.Range("C3").Value = 76.34
$ws.Range("D3").Value = 73.64

$ws.Range("A4").Value = "2020年"
$ws.Range("B4").Value = 80.88
$ws.Range("C4").Value = 77.93000000000001
$ws.Range("D4").Value = 75.37

# Remove the now-obsolete 5th row entirely so the sheet dimension becomes A1:D4.
$ws.Rows.Item(5).Delete()
